$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B33 was stored as text "3" - convert it to a real number (3)
$ws.Cells.Item(33, 2).Value = 3

# New row 34 with Ying Tang's additional annotation
$ws.Cells.Item(34, 1).Value = "Ying Tang"

# B34 keeps the numeric-looking value "1" stored as text (not a number)
$ws.Cells.Item(34, 2).Value = "'1"
$ws.Cells.Item(34, 2).Style = "Normal"

$ws.Cells.Item(34, 3).Value = "NOT a proper"
$ws.Cells.Item(34, 4).Value = "CRT"
$ws.Cells.Item(34, 5).Value = "MET"
$ws.Cells.Item(34, 6).Value = "aa721c36-81b2-451c-915e-fe15286fe992"
$ws.Cells.Item(34, 7).Value = "SygwwGbRW_annotated.xlsx"
$ws.Cells.Item(34, 8).Value = "This is NOT a proper navigation agent."
